$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.019.79'
$ws.Range('D3').Value = '1.420.67'
$ws.Range('E3').Value = '  -7.70%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = "'274.34"
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('D7').Value = "'0.3723"
$ws.Range('E7').Value = '  -3.94%  '
$ws.Range('D8').Value = "'0.3086"
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = "'39.92"
$ws.Range('E9').Value = '  -7.37%  '
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('D11').Value = "'0.06601"
$ws.Range('E11').Value = '  -8.29%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  -3.90%  '
$ws.Range('D14').Value = "'17.16"
$ws.Range('E14').Value = '  -7.60%  '
$ws.Range('D15').Value = "'6.183"
$ws.Range('E15').Value = '  -6.27%  '
$ws.Range('D16').Value = '1.421.03'
$ws.Range('E16').Value = '  -7.68%  '
$ws.Range('D17').Value = "'0.00001009"
$ws.Range('E17').Value = '  -9.15%  '
$ws.Range('D18').Value = "'0.05802"
$ws.Range('E18').Value = '  -11.93%  '
$ws.Range('D19').Value = "'74.62"
$ws.Range('E19').Value = '  -10.34%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = "'5.639"
$ws.Range('D22').Value = "'14.52"
$ws.Range('E22').Value = '  -5.63%  '
$ws.Range('D23').Value = "'11.01"
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').Value = "'2.330"
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').Value = '20.024.56'
$ws.Range('E25').Value = '  -7.95%  '
$ws.Range('D26').Value = "'2.296"
$ws.Range('E26').Value = '  -3.21%  '
$ws.Range('D27').Value = "'139.17"
$ws.Range('E27').Value = '  -5.08%  '
$ws.Range('D28').Value = "'16.91"
$ws.Range('E28').Value = '  -7.89%  '
$ws.Range('D29').Value = '1.580.75'
$ws.Range('E29').Value = '  -7.70%  '
$ws.Range('D30').Value = "'109.25"
$ws.Range('E30').Value = '  -6.92%  '
$ws.Range('D31').Value = "'3.811"
$ws.Range('E31').Value = '  -21.22%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'0.8909"
$ws.Range('E32').Value = '  -8.04%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'5.417"
$ws.Range('E33').Value = '  -8.29%  '
$ws.Range('D34').Value = "'0.07745"
$ws.Range('E34').Value = '  -5.20%  '
$ws.Range('D35').Value = "'8.440"
$ws.Range('E35').Value = '  -4.98%  '
$ws.Range('E36').Value = '  +6.11%  '
$ws.Range('D37').Value = "'0.05751"
$ws.Range('E37').Value = '  -5.07%  '
$ws.Range('D38').Value = "'4.798"
$ws.Range('E38').Value = '  -6.73%  '
$ws.Range('D39').Value = "'0.9999"
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = "'0.1929"
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('D41').Value = "'0.02040"
$ws.Range('E41').Value = '  -7.40%  '
$ws.Range('D42').Value = "'1.083"
$ws.Range('E42').Value = '  -9.08%  '
$ws.Range('D43').Value = "'1.274"
$ws.Range('E43').Value = '  -14.42%  '
$ws.Range('D44').Value = "'0.5330"
$ws.Range('E44').Value = '  -7.17%  '
$ws.Range('D45').Value = "'3.539"
$ws.Range('E45').Value = '  -5.35%  '
$ws.Range('D46').Value = "'12.27"
$ws.Range('E46').Value = '  -5.74%  '
$ws.Range('D47').Value = "'0.5138"
$ws.Range('E47').Value = '  -6.80%  '
$ws.Range('D48').Value = "'1.801"
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('D49').Value = "'109.83"
$ws.Range('E49').Value = '  -6.30%  '
$ws.Range('D50').Value = "'1.052"
$ws.Range('E50').Value = '  -8.00%  '
$ws.Range('D51').Value = "'1.000"
$ws.Range('E51').Value = '  -0.07%  '
